# The "H 72" record (original row 2) is removed from the missing_data sheet.
# Deleting the entire row shifts every following row up by one, which is
# exactly what the target workbook shows (dimension shrinks from F63 to F62
# and every subsequent record's row number decreases by 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
